$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "Tempo 2" timing values (H6:H8)
$ws.Range("H6").Value = 1730.3349800000001
$ws.Range("H7").Value = 1690.43894
$ws.Range("H8").Value = 1677.2099900000001

# Add a new (empty) underlined cell at M21, which extends the used range
# and introduces a new font/cell style in the workbook.
$ws.Range("M21").Font.Underline = 1

# Move the active selection to the newly touched cell.
$ws.Range("M21").Select()
